$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - set value then copy the formatting from the existing
# header cell G1 (bold/border/centered style) so it matches the other
# header cells exactly.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells H2:H7 - new "Save" column values
$saveValues = @(1, 1, 0, 0, 1, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
